$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")

# Convert the chart series' cached name references (currently pointing at the
# now-missing external workbook, e.g. '[1]Template'!$A$3) into literal text,
# exactly like Excel does when you break the link to the source workbook.
$chartObj = $ws.ChartObjects(1)
$chart = $chartObj.Chart
$chart.SeriesCollection(1).Name = "Fully Up-to-Date"
$chart.SeriesCollection(2).Name = "Requires Update"
$chart.SeriesCollection(3).Name = "Total listings for criterion"
$chart.SeriesCollection(4).Name = "Code Sets Up-to-Date"
$chart.SeriesCollection(5).Name = "Functionality Up-to-Date"
$chart.SeriesCollection(6).Name = "Standards Up-to-Date"

# Break the external link itself (removes xl/externalLinks/*, the
# <externalReferences> block in workbook.xml, and flips updateLinks so Excel
# never again tries to refresh it).
$sources = $wb.LinkSources(1)
if ($sources) {
    foreach ($src in $sources) {
        $wb.BreakLink($src, 1)
    }
}

# Match the saved selection (cell A7 on the Template sheet).
$ws.Range("A7").Select()
